# daily auto push: 2026-01-10 22:32 UTC
#
# The log sheet records one row per (date, hour) observation. A new
# observation for 2026/01/11 at hour 5 (ranking 150) needs to be inserted
# right after the existing 2026/01/11 row (row 625), pushing every
# subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 625; everything currently at/after 625 (the
# 2026/12/29 ... 2027/01/05 block) shifts down to 626...667.
$ws.Rows.Item(625).Insert()

# Force column A to be read as literal text (not auto-parsed into a date
# serial) while we fill in the new row's values, matching how the rest of
# the sheet stores its "日付" column as plain text.
$ws.Range("A625").NumberFormat = "@"

$ws.Range("A625").Value = "2026/01/11"
$ws.Range("B625").Value = "日"
$ws.Range("C625").Value = 5
$ws.Range("D625").Value = 150

# Drop back to the default "Normal" style so the new cells don't carry a
# stray text-format style that the original data rows never had.
$ws.Range("A625").Style = "Normal"
